$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Sentinel Defined" column (H) from "No" to "Yes" for these rows
$rows = @(12, 16, 19, 21, 23, 24, 25, 26)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = "Yes"
}

# Update the frozen pane / selection on the sheet view
$ws.Activate()
$ws.Range("H25").Select()
$excel.ActiveWindow.ScrollRow = 22
